# Week 13 logging update
# - Rushing sheet: bump a few players' rushing attempt stats
# - Receiving sheet: bump several players' receiving stats and add a new
#   player (D.Parker) as a row inserted right before J.Waddle

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Rushing sheet
# ----------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# T.Tagovailoa (row 2): 1DATT 3 -> 4
$rushing.Cells.Item(2, 3).Value = 4

# M.Gaskin (row 4): 1DATT 78->87, 2DATT 49->54, 3DATT 11->12, RZATT 25->27
$rushing.Cells.Item(4, 3).Value = 87
$rushing.Cells.Item(4, 4).Value = 54
$rushing.Cells.Item(4, 5).Value = 12
$rushing.Cells.Item(4, 6).Value = 27

# S.Ahmed (row 6): 1DATT 21->24, 2DATT 20->25
$rushing.Cells.Item(6, 3).Value = 24
$rushing.Cells.Item(6, 4).Value = 25

# ----------------------------------------------------------------------
# Receiving sheet
# ----------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# M.Brown (row 2): Short Target 52->54, Short Comp 41->43, RZ Target 7->8, RZ Comp 5->6
$receiving.Cells.Item(2, 3).Value = 54
$receiving.Cells.Item(2, 4).Value = 43
$receiving.Cells.Item(2, 7).Value = 8
$receiving.Cells.Item(2, 8).Value = 6

# Insert a new row for D.Parker right before J.Waddle (currently row 6)
$receiving.Rows("6:6").Insert()

$receiving.Cells.Item(6, 1).Value = 4
$receiving.Cells.Item(6, 1).Style = $receiving.Cells.Item(7, 1).Style
$receiving.Cells.Item(6, 2).Value = "D.Parker"
$receiving.Cells.Item(6, 3).Value = 5
$receiving.Cells.Item(6, 4).Value = 5
$receiving.Cells.Item(6, 5).Value = 0
$receiving.Cells.Item(6, 6).Value = 0
$receiving.Cells.Item(6, 7).Value = 0
$receiving.Cells.Item(6, 8).Value = 0

# J.Waddle (now row 7): Short Target 90->100, Short Comp 72->80, Deep Target 15->16,
#   Deep Comp 5->6, RZ Target 12->13, RZ Comp 9->10
$receiving.Cells.Item(7, 3).Value = 100
$receiving.Cells.Item(7, 4).Value = 80
$receiving.Cells.Item(7, 5).Value = 16
$receiving.Cells.Item(7, 6).Value = 6
$receiving.Cells.Item(7, 7).Value = 13
$receiving.Cells.Item(7, 8).Value = 10

# J.Grant (now row 8): unchanged
# P.Williams (now row 9): unchanged

# A.Wilson (now row 10): Short Target 27->35, Short Comp 19->23
$receiving.Cells.Item(10, 3).Value = 35
$receiving.Cells.Item(10, 4).Value = 23

# M.Hollins (now row 11): Short Target 15->16, Short Comp 9->10, RZ Target 5->6, RZ Comp 3->4
$receiving.Cells.Item(11, 3).Value = 16
$receiving.Cells.Item(11, 4).Value = 10
$receiving.Cells.Item(11, 7).Value = 6
$receiving.Cells.Item(11, 8).Value = 4

# I.Ford (now row 12): Short Target 7->8, Short Comp 5->6, RZ Target 2->3, RZ Comp 1->2
$receiving.Cells.Item(12, 3).Value = 8
$receiving.Cells.Item(12, 4).Value = 6
$receiving.Cells.Item(12, 7).Value = 3
$receiving.Cells.Item(12, 8).Value = 2

# K.Merritt (now row 13): unchanged

# M.Gesicki (now row 14): Short Target 63->72, Short Comp 40->47, Deep Target 15->17,
#   RZ Target 5->7, RZ Comp 5->6
$receiving.Cells.Item(14, 3).Value = 72
$receiving.Cells.Item(14, 4).Value = 47
$receiving.Cells.Item(14, 5).Value = 17
$receiving.Cells.Item(14, 7).Value = 7
$receiving.Cells.Item(14, 8).Value = 6

# A.Shaheen (now row 15): unchanged

# H.Long (now row 16): Short Target 1->2, Short Comp 0->1
$receiving.Cells.Item(16, 3).Value = 2
$receiving.Cells.Item(16, 4).Value = 1

# D.Smythe (now row 17): Short Target 25->26
$receiving.Cells.Item(17, 3).Value = 26
